$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1229557.4
$ws.Range("J17").Value = 1299763.9
$ws.Range("L17").Value = 3899291.7
$ws.Range("N17").Value = -3899627.7
$ws.Range("H80").Value = 96140.336
$ws.Range("I80").Value = 125972.625
$ws.Range("J80").Value = 677
$ws.Range("K80").Value = 377917.875
$ws.Range("L80").Value = 2031
$ws.Range("M80").Value = -376919.875
$ws.Range("N80").Value = -4027
$ws.Range("H83").Value = 96140.336
$ws.Range("I83").Value = 125972.625
$ws.Range("J83").Value = 677
$ws.Range("K83").Value = 1133753.625
$ws.Range("L83").Value = 6093
$ws.Range("M83").Value = -1128761.625
$ws.Range("N83").Value = -16077
$ws.Range("H86").Value = 19690.385
$ws.Range("I86").Value = 2397.4443
$ws.Range("J86").Value = 58599.5
$ws.Range("K86").Value = 2397.4443
$ws.Range("L86").Value = 58599.5
$ws.Range("M86").Value = -1274.4443
$ws.Range("N86").Value = -60845.5
$ws.Range("H88").Value = 8657.154
$ws.Range("I88").Value = 7388.5
$ws.Range("J88").Value = 8887.817999999999
$ws.Range("K88").Value = 7388.5
$ws.Range("L88").Value = 8887.817999999999
$ws.Range("M88").Value = -6982.5
$ws.Range("N88").Value = -9699.817999999999
$ws.Range("H89").Value = 19690.385
$ws.Range("I89").Value = 2397.4443
$ws.Range("J89").Value = 58599.5
$ws.Range("K89").Value = 11987.2215
$ws.Range("L89").Value = 292997.5
$ws.Range("M89").Value = -6371.2215
$ws.Range("N89").Value = -304229.5
$ws.Range("H91").Value = 8657.154
$ws.Range("I91").Value = 7388.5
$ws.Range("J91").Value = 8887.817999999999
$ws.Range("K91").Value = 7388.5
$ws.Range("L91").Value = 8887.817999999999
$ws.Range("M91").Value = -5984.5
$ws.Range("N91").Value = -11695.818
$ws.Range("H116").Value = 2358684.2
$ws.Range("I116").Value = 3533277.8
$ws.Range("J116").Value = 9497.5
$ws.Range("K116").Value = 3533277.8
$ws.Range("L116").Value = 9497.5
$ws.Range("M116").Value = -3529835.8
$ws.Range("N116").Value = -16381.5
$ws.Range("H129").Value = 71429910
$ws.Range("I129").Value = 1198.2727
$ws.Range("J129").Value = 333335200
$ws.Range("K129").Value = 3594.8181
$ws.Range("L129").Value = 1000005600
$ws.Range("M129").Value = 1405.1819
$ws.Range("N129").Value = -1000015600
$ws.Range("H132").Value = 2904.2559
$ws.Range("I132").Value = 2882.4614
$ws.Range("K132").Value = 8647.3842
$ws.Range("M132").Value = -6117.3842
$ws.Range("H137").Value = 8025.1816
$ws.Range("I137").Value = 8654.861999999999
$ws.Range("K137").Value = 25964.586
$ws.Range("M137").Value = -23414.586
$ws.Range("H140").Value = 54000
$ws.Range("J140").Value = 78000
$ws.Range("L140").Value = 78000
$ws.Range("N140").Value = -88360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2106.2856
$ws.Range("I32").Value = 2157.3076
$ws.Range("J32").Value = 1443
$ws.Range("K32").Value = 2157.3076
$ws.Range("L32").Value = 1443
$ws.Range("M32").Value = -1870.3076
$ws.Range("N32").Value = -2017
$ws.Range("H41").Value = 7790
$ws.Range("I41").Value = 5730
$ws.Range("K41").Value = 5730
$ws.Range("M41").Value = -5316
$ws.Range("H45").Value = 99463.09
$ws.Range("I45").Value = 146249.78
$ws.Range("K45").Value = 146249.78
$ws.Range("M45").Value = -145872.78
$ws.Range("H63").Value = 5302.6665
$ws.Range("I63").Value = 5302.6665
$ws.Range("K63").Value = 5302.6665
$ws.Range("M63").Value = -4616.6665
$ws.Range("H66").Value = 5302.6665
$ws.Range("I66").Value = 5302.6665
$ws.Range("K66").Value = 26513.3325
$ws.Range("M66").Value = -23081.3325
$ws.Range("H110").Value = 2409.0667
$ws.Range("I110").Value = 1568.091
$ws.Range("K110").Value = 1568.091
$ws.Range("M110").Value = 476.9090000000001
$ws.Range("H132").Value = 1814.6904
$ws.Range("I132").Value = 961.9143
$ws.Range("J132").Value = 6078.5713
$ws.Range("K132").Value = 2885.7429
$ws.Range("L132").Value = 18235.7139
$ws.Range("M132").Value = -355.7429000000002
$ws.Range("N132").Value = -23295.7139

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6839
$ws.Range("I20").Value = 3328
$ws.Range("K20").Value = 3328
$ws.Range("M20").Value = -3081
$ws.Range("H55").Value = 42000
$ws.Range("J55").Value = 42000
$ws.Range("L55").Value = 42000
$ws.Range("N55").Value = -42546
$ws.Range("H105").Value = 58774.555
$ws.Range("I105").Value = 73138.86
$ws.Range("K105").Value = 73138.86
$ws.Range("M105").Value = -71391.86
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2661.2068
$ws.Range("I31").Value = 1946.7826
$ws.Range("J31").Value = 5399.8335
$ws.Range("K31").Value = 1946.7826
$ws.Range("L31").Value = 5399.8335
$ws.Range("M31").Value = -1651.7826
$ws.Range("N31").Value = -5989.8335
$ws.Range("H34").Value = 2661.2068
$ws.Range("I34").Value = 1946.7826
$ws.Range("J34").Value = 5399.8335
$ws.Range("K34").Value = 1946.7826
$ws.Range("L34").Value = 5399.8335
$ws.Range("M34").Value = -1744.7826
$ws.Range("N34").Value = -5803.8335
$ws.Range("H105").Value = 176283
$ws.Range("I105").Value = 420479.8
$ws.Range("K105").Value = 420479.8
$ws.Range("M105").Value = -418732.8
$ws.Range("H107").Value = 14713.3125
$ws.Range("J107").Value = 1531.625
$ws.Range("L107").Value = 1531.625
$ws.Range("N107").Value = -5371.625
$ws.Range("H121").Value = 54950
$ws.Range("J121").Value = 54950
$ws.Range("L121").Value = 54950
$ws.Range("N121").Value = -57570
$ws.Range("H123").Value = 78000
$ws.Range("J123").Value = 78000
$ws.Range("L123").Value = 78000
$ws.Range("N123").Value = -87800
$ws.Range("H133").Value = 57545
$ws.Range("J133").Value = 57545
$ws.Range("L133").Value = 57545
$ws.Range("N133").Value = -62605

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 324679.25
$ws.Range("I5").Value = 1923
$ws.Range("J5").Value = 627263.25
$ws.Range("K5").Value = 5769
$ws.Range("L5").Value = 1881789.75
$ws.Range("M5").Value = -5657
$ws.Range("N5").Value = -1882013.75
$ws.Range("H81").Value = 10526.083
$ws.Range("I81").Value = 2073.25
$ws.Range("J81").Value = 14752.5
$ws.Range("K81").Value = 6219.75
$ws.Range("L81").Value = 44257.5
$ws.Range("M81").Value = -5096.75
$ws.Range("N81").Value = -46503.5
$ws.Range("H84").Value = 10526.083
$ws.Range("I84").Value = 2073.25
$ws.Range("J84").Value = 14752.5
$ws.Range("K84").Value = 18659.25
$ws.Range("L84").Value = 132772.5
$ws.Range("M84").Value = -13043.25
$ws.Range("N84").Value = -144004.5
$ws.Range("H96").Value = 3999.5
$ws.Range("J96").Value = 3999.5
$ws.Range("L96").Value = 11998.5
$ws.Range("N96").Value = -16116.5
$ws.Range("H98").Value = 945.5833
$ws.Range("I98").Value = 1050.4286
$ws.Range("K98").Value = 3151.2858
$ws.Range("M98").Value = -1653.2858
$ws.Range("H131").Value = 58825056
$ws.Range("I131").Value = 111111830
$ws.Range("K131").Value = 333335490
$ws.Range("M131").Value = -333330450
$ws.Range("H132").Value = 38448.43
$ws.Range("I132").Value = 1071
$ws.Range("K132").Value = 9639
$ws.Range("M132").Value = -7109
$ws.Range("H135").Value = 324679.25
$ws.Range("I135").Value = 1923
$ws.Range("J135").Value = 627263.25
$ws.Range("K135").Value = 17307
$ws.Range("L135").Value = 5645369.25
$ws.Range("M135").Value = -14772
$ws.Range("N135").Value = -5650439.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 42573216
$ws.Range("I11").Value = 53181820
$ws.Range("K11").Value = 53181820
$ws.Range("M11").Value = -53181681
$ws.Range("H63").Value = 68750
$ws.Range("J63").Value = 68750
$ws.Range("L63").Value = 68750
$ws.Range("N63").Value = -70122
$ws.Range("H66").Value = 68750
$ws.Range("J66").Value = 68750
$ws.Range("L66").Value = 206250
$ws.Range("N66").Value = -213114
$ws.Range("H80").Value = 9725.571
$ws.Range("I80").Value = 12615.8
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 12615.8
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -11617.8
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 9725.571
$ws.Range("I83").Value = 12615.8
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 63079
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -58087
$ws.Range("N83").Value = -22484
$ws.Range("H132").Value = 2917.1025
$ws.Range("I132").Value = 2601.9714
$ws.Range("K132").Value = 7805.914199999999
$ws.Range("M132").Value = -5275.914199999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 4050.625
$ws.Range("I20").Value = 3234.1667
$ws.Range("K20").Value = 3234.1667
$ws.Range("M20").Value = -3008.1667
$ws.Range("H40").Value = 30511.666
$ws.Range("I40").Value = 46983.777
$ws.Range("J40").Value = 14039.556
$ws.Range("K40").Value = 46983.777
$ws.Range("L40").Value = 14039.556
$ws.Range("M40").Value = -46847.777
$ws.Range("N40").Value = -14311.556

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H97").Value = 49066
$ws.Range("J97").Value = 49066
$ws.Range("L97").Value = 49066
$ws.Range("N97").Value = -51048
$ws.Range("H132").Value = 10053.128
$ws.Range("I132").Value = 10551.077
$ws.Range("K132").Value = 31653.231
$ws.Range("M132").Value = -29123.231
$ws.Range("H135").Value = 124908.4
$ws.Range("J135").Value = 124908.4
$ws.Range("L135").Value = 124908.4
$ws.Range("N135").Value = -135048.4
